$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Vehicle category code mapping rules (rows 2-9), header already in row 1.
$data = @(
    @("载货汽车", 1),
    @("越野汽车", 2),
    @("自卸汽车", 3),
    @("牵引汽车", 4),
    @("专用汽车", 5),
    @("客车", 6),
    @("轿车", 7),
    @("半挂车", 9)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r++
}

# Update selection to reflect the next empty row, as seen in the target workbook.
$ws.Range("B10").Select()
